# Updated cryptos list: Price (D) and Volume(1h) (E) columns
# D-column values are plain numeric-looking text in the source data, so we
# force a Text number format before assigning then clear the format again so
# the value is preserved as a literal string without leaving a style behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.203.27"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.671.43"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.59%  "

$ws.Range("E4").Value = "  -0.62%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.43"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5278"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.71%  "

$ws.Range("E7").Value = "  -0.59%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2642"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06281"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.31"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07554"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.676.92"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.453"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5596"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.10"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008000"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.239.64"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.88%  "

$ws.Range("E18").Value = "  -0.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.790"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.28"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.40"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.196"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.73"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1257"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.574"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.94"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06178"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("E29").Value = "  -1.88%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.284"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.501"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.13%  "

$ws.Range("E32").Value = "  -4.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.629"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9993"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6069"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.412"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.743"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.139"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01621"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.098.30"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8767"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.59%  "

$ws.Range("E42").Value = "  -0.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.80"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.822.90"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000112"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.93"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.09%  "

$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.032"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05228"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4254"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.983"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.37%  "

